# "Modificacion en etiqueta por acento" — fix a mislabeled/incomplete status
# column on the "Hoja1" task sheet:
#   - Row 17 ("Factura A y B") was missing its Responsable/Estado entries.
#   - Row 31 ("Agregar patron fechas...") was missing its Estado entry.
#   - Row 40 ("acentos!") was missing its Responsable/Porcentaje entries.
# Also restores the sheet view to its default top-left position with the
# selection resting on A19 instead of the stray B34:B35 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 17: Responsable = Agustina, Estado = en proceso
$ws.Range("B17").Value = "Agustina"
$ws.Range("C17").Value = "en proceso"
# C17 previously held an (empty) underlined-style placeholder cell; drop the
# underline now that it carries real text, matching the rest of the "Estado"
# column (e.g. C4, C12) which use the plain, non-underlined style.
$ws.Range("C17").Font.Underline = $false

# Row 31: Estado = en proceso
$ws.Range("C31").Value = "en proceso"

# Row 40: Responsable = Lucas, Porcentaje = 100%
$ws.Range("B40").Value = "Lucas"
$ws.Range("C40").Value = 1
$ws.Range("C40").NumberFormat = "0%"

# Reset the view: select A19 (instead of the previous topLeftCell="A19" +
# B34:B35 selection) which also drops the stale topLeftCell scroll position.
$ws.Range("A19").Select()
